$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 49 (last row removed)
$ws.Rows("49").Delete()

# Update B and C values for rows 2-48
$ws.Range("B2").Value = 2.431696818617393
$ws.Range("C2").Value = 1.674599064774749
$ws.Range("B3").Value = 3.007614532213715
$ws.Range("C3").Value = 3.984993748082846
$ws.Range("B4").Value = 5.854059708344627
$ws.Range("C4").Value = 6.489268202224746
$ws.Range("B5").Value = 9.052376654460819
$ws.Range("C5").Value = 8.234356864579766
$ws.Range("B6").Value = 10.65228171639642
$ws.Range("C6").Value = 9.986666524096563
$ws.Range("B7").Value = 12.79092742448272
$ws.Range("C7").Value = 12.14304170838767
$ws.Range("B8").Value = 14.01198498540949
$ws.Range("C8").Value = 14.41969045064456
$ws.Range("B9").Value = 14.94263752795537
$ws.Range("C9").Value = 16.39139445198223
$ws.Range("B10").Value = 15.1992433257655
$ws.Range("C10").Value = 18.05182920597896
$ws.Range("B11").Value = 16.28431624865659
$ws.Range("C11").Value = 20.17413857306362
$ws.Range("B12").Value = 16.49394836035179
$ws.Range("C12").Value = 22.14556670371766
$ws.Range("B13").Value = 19.22949956528216
$ws.Range("C13").Value = 24.10719449608797
$ws.Range("B14").Value = 27.85341600608195
$ws.Range("C14").Value = 26.24593323884858
$ws.Range("B15").Value = 32.73780626499064
$ws.Range("C15").Value = 27.9510554652646
$ws.Range("B16").Value = 34.23868408370404
$ws.Range("C16").Value = 30.27125258116521
$ws.Range("B17").Value = 35.41211055860812
$ws.Range("C17").Value = 32.09976626311879
$ws.Range("B18").Value = 36.69596135919109
$ws.Range("C18").Value = 34.37236182290327
$ws.Range("B19").Value = 38.12030398132588
$ws.Range("C19").Value = 36.01648844162155
$ws.Range("B20").Value = 39.16341522200401
$ws.Range("C20").Value = 38.62341746199056
$ws.Range("B21").Value = 40.71825625400218
$ws.Range("C21").Value = 41.11153218194475
$ws.Range("B22").Value = 44.64775902076117
$ws.Range("C22").Value = 43.33549247313314
$ws.Range("B23").Value = 46.59940782102131
$ws.Range("C23").Value = 45.22238845298661
$ws.Range("B24").Value = 49.84312402721868
$ws.Range("C24").Value = 48.00894816927949
$ws.Range("B25").Value = 50.22276708884628
$ws.Range("C25").Value = 49.86270241209153
$ws.Range("B26").Value = 53.67304836397857
$ws.Range("C26").Value = 52.5684273540824
$ws.Range("B27").Value = 57.40189534648996
$ws.Range("C27").Value = 54.5802741726667
$ws.Range("B28").Value = 58.78797062809473
$ws.Range("C28").Value = 56.42207170184363
$ws.Range("B29").Value = 61.26257758814375
$ws.Range("C29").Value = 59.26344984936472
$ws.Range("B30").Value = 63.5182907454618
$ws.Range("C30").Value = 61.32844623626853
$ws.Range("B31").Value = 68.385957980635
$ws.Range("C31").Value = 63.31042775792989
$ws.Range("B32").Value = 71.38338409566147
$ws.Range("C32").Value = 65.15172631515384
$ws.Range("B33").Value = 72.36999207921843
$ws.Range("C33").Value = 67.04107793064907
$ws.Range("B34").Value = 73.21514856293743
$ws.Range("C34").Value = 68.91723381964916
$ws.Range("B35").Value = 73.98281501822837
$ws.Range("C35").Value = 71.27747138957832
$ws.Range("B36").Value = 76.03595283560945
$ws.Range("C36").Value = 74.17458275838875
$ws.Range("B37").Value = 79.9357122905369
$ws.Range("C37").Value = 76.00120348587149
$ws.Range("B38").Value = 81.41384833472864
$ws.Range("C38").Value = 77.63429882872356
$ws.Range("B39").Value = 82.74767445380279
$ws.Range("C39").Value = 80.23771661791326
$ws.Range("B40").Value = 83.277899877884
$ws.Range("C40").Value = 82.66403401824664
$ws.Range("B41").Value = 83.93942891712523
$ws.Range("C41").Value = 84.69859280291925
$ws.Range("B42").Value = 88.35446693740843
$ws.Range("C42").Value = 86.90955604336531
$ws.Range("B43").Value = 89.98424039777574
$ws.Range("C43").Value = 88.93354753075404
$ws.Range("B44").Value = 90.96418754052027
$ws.Range("C44").Value = 90.84001761490315
$ws.Range("B45").Value = 94.03486452489147
$ws.Range("C45").Value = 93.54530125549694
$ws.Range("B46").Value = 94.97025942447965
$ws.Range("C46").Value = 95.605371126582
$ws.Range("B47").Value = 95.13396070046596
$ws.Range("C47").Value = 97.49886596242141
$ws.Range("B48").Value = 99.4846855292158
$ws.Range("C48").Value = 99.55988490752101
